$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue ($ws.Cells.Item(2, 4)) '42.762.37'
Set-TextValue ($ws.Cells.Item(2, 5)) '  -1.08%  '

# Row 3
Set-TextValue ($ws.Cells.Item(3, 4)) '2.540.40'
Set-TextValue ($ws.Cells.Item(3, 5)) '  -1.64%  '

# Row 4
Set-TextValue ($ws.Cells.Item(4, 5)) '  +0.04%  '

# Row 5
Set-TextValue ($ws.Cells.Item(5, 4)) '308.58'
Set-TextValue ($ws.Cells.Item(5, 5)) '  -2.34%  '

# Row 6
Set-TextValue ($ws.Cells.Item(6, 4)) '100.86'
Set-TextValue ($ws.Cells.Item(6, 5)) '  +4.17%  '

# Row 7
Set-TextValue ($ws.Cells.Item(7, 4)) '0.572'
Set-TextValue ($ws.Cells.Item(7, 5)) '  -0.95%  '

# Row 8
Set-TextValue ($ws.Cells.Item(8, 5)) '  +0.17%  '

# Row 9
Set-TextValue ($ws.Cells.Item(9, 5)) '  -2.09%  '

# Row 10
Set-TextValue ($ws.Cells.Item(10, 4)) '36.30'
Set-TextValue ($ws.Cells.Item(10, 5)) '  +1.89%  '

# Row 11
Set-TextValue ($ws.Cells.Item(11, 4)) '0.0805'
Set-TextValue ($ws.Cells.Item(11, 5)) '  -1.13%  '

# Row 12
Set-TextValue ($ws.Cells.Item(12, 4)) '7.36'
Set-TextValue ($ws.Cells.Item(12, 5)) '  -1.79%  '

# Row 13
Set-TextValue ($ws.Cells.Item(13, 5)) '  -0.10%  '

# Row 14
Set-TextValue ($ws.Cells.Item(14, 4)) '2.923.11'
Set-TextValue ($ws.Cells.Item(14, 5)) '  -1.84%  '

# Row 15
Set-TextValue ($ws.Cells.Item(15, 4)) '15.84'
Set-TextValue ($ws.Cells.Item(15, 5)) '  +4.49%  '

# Row 16
Set-TextValue ($ws.Cells.Item(16, 4)) '2.569.50'
Set-TextValue ($ws.Cells.Item(16, 5)) '  +2.74%  '

# Row 17
Set-TextValue ($ws.Cells.Item(17, 4)) '0.813'
Set-TextValue ($ws.Cells.Item(17, 5)) '  -3.72%  '

# Row 18
Set-TextValue ($ws.Cells.Item(18, 4)) '42.740.21'
Set-TextValue ($ws.Cells.Item(18, 5)) '  -1.23%  '

# Row 19
Set-TextValue ($ws.Cells.Item(19, 5)) '  -0.91%  '

# Row 20
Set-TextValue ($ws.Cells.Item(20, 5)) '  -0.97%  '

# Row 21
Set-TextValue ($ws.Cells.Item(21, 4)) '12.27'
Set-TextValue ($ws.Cells.Item(21, 5)) '  -2.23%  '

# Row 22
Set-TextValue ($ws.Cells.Item(22, 4)) '69.36'
Set-TextValue ($ws.Cells.Item(22, 5)) '  -0.17%  '

# Row 23
Set-TextValue ($ws.Cells.Item(23, 4)) '244.37'

# Row 24
Set-TextValue ($ws.Cells.Item(24, 5)) '  -2.72%  '

# Row 25
Set-TextValue ($ws.Cells.Item(25, 5)) '  -1.34%  '

# Row 26
Set-TextValue ($ws.Cells.Item(26, 2)) 'EthereumClassic'
Set-TextValue ($ws.Cells.Item(26, 3)) 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue ($ws.Cells.Item(26, 4)) '26.12'
Set-TextValue ($ws.Cells.Item(26, 5)) '  -4.29%  '

# Row 27
Set-TextValue ($ws.Cells.Item(27, 2)) 'Dai'
Set-TextValue ($ws.Cells.Item(27, 3)) 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue ($ws.Cells.Item(27, 4)) '0.936'
Set-TextValue ($ws.Cells.Item(27, 5)) '  -6.41%  '

# Row 28
Set-TextValue ($ws.Cells.Item(28, 5)) '  -5.15%  '

# Row 29
Set-TextValue ($ws.Cells.Item(29, 4)) '39.24'
Set-TextValue ($ws.Cells.Item(29, 5)) '  -2.36%  '

# Row 30
Set-TextValue ($ws.Cells.Item(30, 4)) '10.20'
Set-TextValue ($ws.Cells.Item(30, 5)) '  -1.46%  '

# Row 31
Set-TextValue ($ws.Cells.Item(31, 4)) '5.80'
Set-TextValue ($ws.Cells.Item(31, 5)) '  -0.86%  '

# Row 32
Set-TextValue ($ws.Cells.Item(32, 4)) '155.80'
Set-TextValue ($ws.Cells.Item(32, 5)) '  +0.52%  '

# Row 33
Set-TextValue ($ws.Cells.Item(33, 5)) '  +12.84%  '

# Row 34
Set-TextValue ($ws.Cells.Item(34, 4)) '0.0793'
Set-TextValue ($ws.Cells.Item(34, 5)) '  -1.55%  '

# Row 35
Set-TextValue ($ws.Cells.Item(35, 5)) '  -2.74%  '

# Row 36
Set-TextValue ($ws.Cells.Item(36, 5)) '  -5.45%  '

# Row 37
Set-TextValue ($ws.Cells.Item(37, 4)) '18.43'
Set-TextValue ($ws.Cells.Item(37, 5)) '  -1.57%  '

# Row 38
Set-TextValue ($ws.Cells.Item(38, 4)) '3.18'
Set-TextValue ($ws.Cells.Item(38, 5)) '  -6.88%  '

# Row 39
Set-TextValue ($ws.Cells.Item(39, 5)) '  +0.25%  '

# Row 40
Set-TextValue ($ws.Cells.Item(40, 5)) '  +0.54%  '

# Row 41
Set-TextValue ($ws.Cells.Item(41, 4)) '4.33'
Set-TextValue ($ws.Cells.Item(41, 5)) '  +9.51%  '

# Row 42
Set-TextValue ($ws.Cells.Item(42, 4)) '22.25'
Set-TextValue ($ws.Cells.Item(42, 5)) '  -0.99%  '

# Row 43
Set-TextValue ($ws.Cells.Item(43, 5)) '  +0.01%  '

# Row 44
Set-TextValue ($ws.Cells.Item(44, 4)) '3.31'
Set-TextValue ($ws.Cells.Item(44, 5)) '  +1.91%  '

# Row 45
Set-TextValue ($ws.Cells.Item(45, 5)) '  -1.66%  '

# Row 46
Set-TextValue ($ws.Cells.Item(46, 4)) '1.958.44'
Set-TextValue ($ws.Cells.Item(46, 5)) '  -2.46%  '

# Row 47
Set-TextValue ($ws.Cells.Item(47, 4)) '8.90'
Set-TextValue ($ws.Cells.Item(47, 5)) '  -0.49%  '

# Row 48
Set-TextValue ($ws.Cells.Item(48, 2)) 'RocketPoolETH'
Set-TextValue ($ws.Cells.Item(48, 3)) 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue ($ws.Cells.Item(48, 4)) '2.764.97'
Set-TextValue ($ws.Cells.Item(48, 5)) '  -2.30%  '

# Row 49
Set-TextValue ($ws.Cells.Item(49, 4)) '0.193'
Set-TextValue ($ws.Cells.Item(49, 5)) '  -0.16%  '

# Row 50
Set-TextValue ($ws.Cells.Item(50, 2)) 'BitcoinSV'
Set-TextValue ($ws.Cells.Item(50, 3)) 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue ($ws.Cells.Item(50, 4)) '80.77'
Set-TextValue ($ws.Cells.Item(50, 5)) '  -2.75%  '

# Row 51
Set-TextValue ($ws.Cells.Item(51, 4)) '0.857'
Set-TextValue ($ws.Cells.Item(51, 5)) '  +10.04%  '
